$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy number-format (date/time) styles down to the new rows (A123:B180) ---
$ws.Range("A2:B2").Copy()
$ws.Range("A123:B180").PasteSpecial(-4122)

# --- Populate the new log rows (123-180) with date/time and event data ---
$ws.Cells.Item(123, 1).Value = 44524
$ws.Cells.Item(123, 2).Value = 0.29166666666666669
$ws.Cells.Item(123, 3).Value = 0
$ws.Cells.Item(124, 1).Value = 44524
$ws.Cells.Item(124, 2).Value = 0.29236111111111113
$ws.Cells.Item(124, 4).Value = 0
$ws.Cells.Item(125, 1).Value = 44524
$ws.Cells.Item(125, 2).Value = 0.31597222222222221
$ws.Cells.Item(125, 5).Value = 1
$ws.Cells.Item(126, 1).Value = 44524
$ws.Cells.Item(126, 2).Value = 0.39583333333333331
$ws.Cells.Item(127, 1).Value = 44524
$ws.Cells.Item(127, 2).Value = 0.47638888888888892
$ws.Cells.Item(128, 1).Value = 44524
$ws.Cells.Item(128, 2).Value = 0.48541666666666666
$ws.Cells.Item(128, 4).Value = 1
$ws.Cells.Item(129, 1).Value = 44524
$ws.Cells.Item(129, 2).Value = 0.4861111111111111
$ws.Cells.Item(129, 5).Value = 1
$ws.Cells.Item(130, 1).Value = 44524
$ws.Cells.Item(130, 2).Value = 0.60763888888888895
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(131, 1).Value = 44524
$ws.Cells.Item(131, 2).Value = 0.61319444444444449
$ws.Cells.Item(131, 3).Value = 0
$ws.Cells.Item(132, 1).Value = 44524
$ws.Cells.Item(132, 2).Value = 0.71527777777777779
$ws.Cells.Item(132, 3).Value = 0
$ws.Cells.Item(133, 1).Value = 44524
$ws.Cells.Item(133, 2).Value = 0.77083333333333337
$ws.Cells.Item(133, 3).Value = 1
$ws.Cells.Item(134, 1).Value = 44524
$ws.Cells.Item(134, 2).Value = 0.79166666666666663
$ws.Cells.Item(134, 5).Value = 1
$ws.Cells.Item(135, 1).Value = 44524
$ws.Cells.Item(135, 2).Value = 0.79513888888888884
$ws.Cells.Item(135, 3).Value = 0
$ws.Cells.Item(136, 1).Value = 44524
$ws.Cells.Item(136, 2).Value = 0.79861111111111116
$ws.Cells.Item(136, 4).Value = 1
$ws.Cells.Item(137, 1).Value = 44524
$ws.Cells.Item(137, 2).Value = 0.84791666666666676
$ws.Cells.Item(137, 3).Value = 1
$ws.Cells.Item(138, 1).Value = 44524
$ws.Cells.Item(138, 2).Value = 0.85138888888888886
$ws.Cells.Item(138, 4).Value = 1
$ws.Cells.Item(139, 1).Value = 44525
$ws.Cells.Item(139, 2).Value = 0.24305555555555555
$ws.Cells.Item(139, 3).Value = 1
$ws.Cells.Item(140, 1).Value = 44525
$ws.Cells.Item(140, 2).Value = 0.3
$ws.Cells.Item(140, 5).Value = 1
$ws.Cells.Item(141, 1).Value = 44525
$ws.Cells.Item(141, 2).Value = 0.31527777777777777
$ws.Cells.Item(141, 4).Value = 1
$ws.Cells.Item(142, 1).Value = 44525
$ws.Cells.Item(142, 2).Value = 0.45833333333333331
$ws.Cells.Item(142, 3).Value = 1
$ws.Cells.Item(143, 1).Value = 44525
$ws.Cells.Item(143, 2).Value = 0.47222222222222227
$ws.Cells.Item(143, 5).Value = 1
$ws.Cells.Item(144, 1).Value = 44525
$ws.Cells.Item(144, 2).Value = 0.4826388888888889
$ws.Cells.Item(144, 3).Value = 1
$ws.Cells.Item(145, 1).Value = 44525
$ws.Cells.Item(145, 2).Value = 0.48333333333333334
$ws.Cells.Item(145, 4).Value = 1
$ws.Cells.Item(146, 1).Value = 44525
$ws.Cells.Item(146, 2).Value = 0.60763888888888895
$ws.Cells.Item(146, 3).Value = 1
$ws.Cells.Item(147, 1).Value = 44525
$ws.Cells.Item(147, 2).Value = 0.76388888888888884
$ws.Cells.Item(147, 3).Value = 1
$ws.Cells.Item(148, 1).Value = 44525
$ws.Cells.Item(148, 2).Value = 0.79166666666666663
$ws.Cells.Item(148, 5).Value = 1
$ws.Cells.Item(149, 1).Value = 44525
$ws.Cells.Item(149, 2).Value = 0.8125
$ws.Cells.Item(149, 4).Value = 1
$ws.Cells.Item(150, 1).Value = 44525
$ws.Cells.Item(150, 2).Value = 0.90625
$ws.Cells.Item(150, 3).Value = 0
$ws.Cells.Item(151, 1).Value = 44526
$ws.Cells.Item(151, 2).Value = 0.17361111111111113
$ws.Cells.Item(151, 3).Value = 1
$ws.Cells.Item(152, 1).Value = 44526
$ws.Cells.Item(152, 2).Value = 0.30763888888888891
$ws.Cells.Item(152, 3).Value = 1
$ws.Cells.Item(153, 1).Value = 44526
$ws.Cells.Item(153, 2).Value = 0.3125
$ws.Cells.Item(153, 5).Value = 1
$ws.Cells.Item(154, 1).Value = 44526
$ws.Cells.Item(154, 2).Value = 0.34375
$ws.Cells.Item(154, 3).Value = 1
$ws.Cells.Item(155, 1).Value = 44526
$ws.Cells.Item(155, 2).Value = 0.34722222222222227
$ws.Cells.Item(155, 4).Value = 1
$ws.Cells.Item(156, 1).Value = 44526
$ws.Cells.Item(156, 2).Value = 0.41597222222222219
$ws.Cells.Item(156, 3).Value = 0
$ws.Cells.Item(157, 1).Value = 44526
$ws.Cells.Item(157, 2).Value = 0.50069444444444444
$ws.Cells.Item(157, 3).Value = 0
$ws.Cells.Item(158, 1).Value = 44526
$ws.Cells.Item(158, 2).Value = 0.53888888888888886
$ws.Cells.Item(158, 5).Value = 1
$ws.Cells.Item(159, 1).Value = 44526
$ws.Cells.Item(159, 2).Value = 0.54583333333333328
$ws.Cells.Item(159, 3).Value = 0
$ws.Cells.Item(160, 1).Value = 44526
$ws.Cells.Item(160, 2).Value = 0.5625
$ws.Cells.Item(160, 4).Value = 0
$ws.Cells.Item(161, 1).Value = 44526
$ws.Cells.Item(161, 2).Value = 0.61388888888888882
$ws.Cells.Item(161, 3).Value = 0
$ws.Cells.Item(162, 1).Value = 44526
$ws.Cells.Item(162, 2).Value = 0.67361111111111116
$ws.Cells.Item(162, 3).Value = 1
$ws.Cells.Item(163, 1).Value = 44526
$ws.Cells.Item(163, 2).Value = 0.74444444444444446
$ws.Cells.Item(163, 3).Value = 0
$ws.Cells.Item(164, 1).Value = 44526
$ws.Cells.Item(164, 2).Value = 0.79513888888888884
$ws.Cells.Item(164, 5).Value = 1
$ws.Cells.Item(165, 1).Value = 44526
$ws.Cells.Item(165, 2).Value = 0.80208333333333337
$ws.Cells.Item(165, 3).Value = 0
$ws.Cells.Item(166, 1).Value = 44526
$ws.Cells.Item(166, 2).Value = 0.81180555555555556
$ws.Cells.Item(166, 4).Value = 1
$ws.Cells.Item(167, 1).Value = 44526
$ws.Cells.Item(167, 2).Value = 0.86388888888888893
$ws.Cells.Item(167, 3).Value = 1
$ws.Cells.Item(168, 1).Value = 44526
$ws.Cells.Item(168, 2).Value = 0.89583333333333337
$ws.Cells.Item(168, 3).Value = 0
$ws.Cells.Item(169, 1).Value = 44527
$ws.Cells.Item(169, 2).Value = 0.26874999999999999
$ws.Cells.Item(169, 3).Value = 1
$ws.Cells.Item(170, 1).Value = 44527
$ws.Cells.Item(170, 2).Value = 0.27777777777777779
$ws.Cells.Item(170, 5).Value = 1
$ws.Cells.Item(171, 1).Value = 44527
$ws.Cells.Item(171, 2).Value = 0.28055555555555556
$ws.Cells.Item(171, 3).Value = 1
$ws.Cells.Item(172, 1).Value = 44527
$ws.Cells.Item(172, 2).Value = 0.42430555555555555
$ws.Cells.Item(172, 3).Value = 0
$ws.Cells.Item(173, 1).Value = 44527
$ws.Cells.Item(173, 2).Value = 0.45833333333333331
$ws.Cells.Item(173, 3).Value = 1
$ws.Cells.Item(174, 1).Value = 44527
$ws.Cells.Item(174, 2).Value = 0.56597222222222221
$ws.Cells.Item(174, 3).Value = 1
$ws.Cells.Item(175, 1).Value = 44527
$ws.Cells.Item(175, 2).Value = 0.60416666666666663
$ws.Cells.Item(175, 3).Value = 1
$ws.Cells.Item(176, 1).Value = 44527
$ws.Cells.Item(176, 2).Value = 0.65625
$ws.Cells.Item(176, 3).Value = 1
$ws.Cells.Item(177, 1).Value = 44527
$ws.Cells.Item(177, 2).Value = 0.73541666666666661
$ws.Cells.Item(177, 3).Value = 1
$ws.Cells.Item(178, 1).Value = 44527
$ws.Cells.Item(178, 2).Value = 0.75347222222222221
$ws.Cells.Item(178, 4).Value = 1
$ws.Cells.Item(179, 1).Value = 44527
$ws.Cells.Item(179, 2).Value = 0.77777777777777779
$ws.Cells.Item(179, 5).Value = 1
$ws.Cells.Item(180, 1).Value = 44527
$ws.Cells.Item(180, 2).Value = 0.78472222222222221
$ws.Cells.Item(180, 3).Value = 0

# --- Add the trailing empty styled row at the very bottom of the sheet (A1048576) ---
$ws.Range("A2").Copy()
$ws.Range("A1048576").PasteSpecial(-4122)

# --- Update the frozen-pane / scroll view: freeze the header row and scroll to the
#     most recently added entries (around row 165). ---
$ws.Range("A2").Select()
$window = $excel.ActiveWindow
$window.FreezePanes = $true
$ws.Range("A165").Select()
